$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.767.98'
$ws.Range("E2").Value = '  +0.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.088.50'
$ws.Range("E3").Value = '  +0.62%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.95'
$ws.Range("E5").Value = '  +0.17%  '

$ws.Range("E6").Value = '  +0.57%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '58.27'
$ws.Range("E8").Value = '  +0.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.394'
$ws.Range("E9").Value = '  +1.08%  '

$ws.Range("E10").Value = '  -0.10%  '

$ws.Range("E11").Value = '  +3.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.20'
$ws.Range("E12").Value = '  +2.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.395.76'
$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.20'
$ws.Range("E14").Value = '  +1.56%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.779'
$ws.Range("E15").Value = '  +1.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.35'
$ws.Range("E16").Value = '  +1.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.077.22'
$ws.Range("E17").Value = '  +0.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.765.51'
$ws.Range("E18").Value = '  +0.49%  '

$ws.Range("E19").Value = '  -0.95%  '

$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("E21").Value = '  +0.74%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.79'
$ws.Range("E22").Value = '  +0.99%  '

$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("E24").Value = '  -0.83%  '

$ws.Range("E25").Value = '  +0.71%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.79'
$ws.Range("E26").Value = '  +9.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '171.30'
$ws.Range("E27").Value = '  +1.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.134'
$ws.Range("E28").Value = '  -3.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.55'
$ws.Range("E29").Value = '  +0.69%  '

$ws.Range("E30").Value = '  +0.26%  '

$ws.Range("E31").Value = '  +1.07%  '

$ws.Range("E32").Value = '  +1.03%  '

$ws.Range("E33").Value = '  +1.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.61'
$ws.Range("E34").Value = '  -0.43%  '

$ws.Range("E35").Value = '  +1.22%  '

$ws.Range("E36").Value = '  -0.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.32'
$ws.Range("E37").Value = '  -1.87%  '

$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("E39").Value = '  +0.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0236'
$ws.Range("E40").Value = '  +9.87%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '101.08'
$ws.Range("E41").Value = '  +3.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0971'
$ws.Range("E42").Value = '  -0.47%  '

$ws.Range("E44").Value = '  +1.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.68'
$ws.Range("E45").Value = '  +1.73%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.453.93'
$ws.Range("E46").Value = '  +0.30%  '

$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.11'
$ws.Range("E47").Value = '  -3.47%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.06'
$ws.Range("E48").Value = '  +0.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.22'
$ws.Range("E49").Value = '  -1.92%  '

$ws.Range("E50").Value = '  -2.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.279.20'
$ws.Range("E51").Value = '  +0.44%  '
